# Daily update of covid19 tracker data files
# Updates country-update narratives on the "Country Updates" sheet for
# Chile, Finland, India, Netherlands and Spain with early-April news.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Updates")

$text_H12 = @'
The Government of Chile has announced several measures in response to the pandemic. First, the government has ensured fiscal resources for detection and treatment of the virus. Second, policy measures have been geared towards supporting household incomes, especially for the most vulnerable, and protecting employment. Third, a broad set of measures point towards injecting liquidity to the economy, thereby facilitating access to working capital, especially for SMEs. Passed the 19th of March, those measures have been reinforced the 8th of April.
'@
$ws.Range("H12").Value = $text_H12

$text_I12 = @'
- Deferral of the payment of Property Tax to be due in April, with 0% interest rate. Applicable to individuals owners of real state with fiscal value lower than CL $133.000.000 (aprox. US$ 160.000). The payment delayed shall be paid along with the 3 next installments to be due in June, September and November 2020.
- Flexibility measures for tax debts payments of low income individuals (flexible payment agreements, with no interest nor penalties). 
- A subsidy (cash bonus) for people without formal work will be provided which will benefit 2 million people.
- Accelerated income tax refund for self-employed (in April 2020 instead of May 2020). 
- Advanced income tax refund for self-employed of amounts withheld in January and February 2020 (rate of 10,75%).
- The government is proposing a bill for the protection of labour income, which will allow the temporary suspension of employees´ functions in the company when -in the framework of this emergency and by mandate of the health authority- the employee must remain at home without the possibility of teleworking. In these cases, the contractual relationship and all corresponding employment rights will be maintained (except the salary), but the employee will be paid from the unemployment insurance under more flexible eligibility criteria. As long as this exceptional situation continues, the employer will maintain the obligation to pay the employee's social security and health contributions. 
- The government will use its legal power to accelerate the discussion in Congress of the bill on employment protection, which allows for the reduction of working hours, compensating the decrease in remuneration with resources from the Solidarity Unemployment Fund.-Creation of a new fund (0;8% of GDP) to protect the income of the most vulnerables. This fund wil benefit up to 3 millions of individuals and the amount they will received will depend on persona circumstances 
'@
$ws.Range("I12").Value = $text_I12

$text_J12 = @'
- Suspension of mandatory monthly provisional payments on account of CIT for the next 3 months. This measure should benefit 700.000 businesses. 
- Deferral of VAT payments for the next 3 months. Applicable to businesses with sales of less than UF 350.000 (aprox. US $12 million). They will be able to pay the VAT in 6 to 12 monthly installments (depending on their size) at a 0% interest rate. This measure should benefit 240.000 businesses.
- Accelerated income tax refund for SMEs (in April 2020 instead of May 2020). This measure should benefit 500.000 SMEs.
- Postponment from April to July 2020 of the payment of CIT for SMEs. The CIT return shall be submitted in April (as usual). This measure should benefit 140.000 SMEs.  
- Deferral of the payment of Property Tax to be due in April 2020, with 0% interest rate. Applicable to businesses with sales lower than UF 350.000 (aprox. US $12 million). The payment delayed shall be paid along with the 3 next installments to be due in June, September and November 2020.  
- Deferral of the payment of Property Tax to be due in April 2020, with 0% interest rate. Applicable to businesses with sales lower than UF 350.000 (aprox. US $12 million). The payment delayed shall be paid along with the 3 next installments to be due in June, September and November 2020.-Plan to guarantee access to credit for 99.8% of firms (84% of jobs), operated through the Fondo de Garantía de Pequeños Empresarios (FOGAPE), with an amount mobilized of 1.2% of GDP. Overall, it is now 10% of GDP that is availabe through FOGAPE for firms to avoid liquidity crunch and preserve them.
'@
$ws.Range("J12").Value = $text_J12

$text_D20 = @'
From 19 March, passenger traffic has been restricted at the Finnish external borders, but Finnish citizens are allowed to return home. Those returning should undertake a mandatory two-week quarantine. From 27 March, airport arrivals have been quarantined (details not announced yet). On 8 April, the government recommended that all shipping companies suspend the sales of tickets to passengers travelling to Finland from Germany, Estonai and Sweden between 11 April and 13 May. The recommendation does not apply to cargo transport. 
'@
$ws.Range("D20").Value = $text_D20

$text_G26 = @'
Early March, the government introduced export restrictions on 26 pharma products (active pharma ingredients and formulations, including paracetamol) to keep domestic supplies steady.   

March 19: Hospitals have been asked to defer elective surgeries. The government capped the price of sanitizer bottles and banned exports of masks, ventilator machines and also textiles used to make masks. 

March 24: An emergency health fund of INR 150 billion (USD 2 billion) for treating COVID-19 patients and strengthening the medical infrastructure of the country including for rapidly ramping up the number of testing facilities, Personal Protective Equipment (PPE), Isolation Beds, ICU beds, Ventilators and other essential equipment announced. Simultaneously, training of medical and paramedical manpower will also be undertaken.Cumulatively 187,904 persons are under surveillance, and around 35,073 have completed 28 days observation period as on day. Of the total samples tested (12872), 2023 samples have been tested by NCDC of which, 52 have tested positive for COVID-19. States such as Gujarat, Assam, Jharkhand, Rajasthan, Goa, Karnataka, MP and J&K are setting up hospitals dedicated for the management of COVID-19. 118 laboratories have been included in the ICMR network of COVID-19 testing with a capacity to test 12,000 samples per day. In the last five days (March 20-24), on an average 1338 samples have been tested per day. Furthermore, 22 private lab chains have registered with ICMR (till 24th March 2020) for COVID-19 testing. They have 15,500 collection centers nationwide. 

April 7: India lifted restrictions on exports of 24 drugs amid coronavirus and allowed limited exports of hydroxychloroquine and paracetamol. It also placed restrictions on the export of most diagnostic testing kits.



'@
$ws.Range("G26").Value = $text_G26

$text_J26 = @'
March 24: The deadline for filing returns of goods and services tax for March, April and May has been extended till June 30. For companies with less than INR 50 million turnover, no interest, penalty and late fee to be charged on late GST return filing. For companies with turnover of over INR 50 million, no late fee and penalty will be charged on late GST return filing; interest rate was reduced to 9%. 

March 26: Government to pay Employee Provident Fund contributions of both employee and employer for the next three months, for all those establishments which have up to with 100 employees, 90% of whom earn less than INR 15,000 a month. 
April 8: To grant immediate relief to taxpayers, the central government announced that all pending income tax refunds up to INR 500000 (about USD 6500) & GST/Custom refunds are to be issued immediately.  

'@
$ws.Range("J26").Value = $text_J26

$text_I38 = @'
- Temporary compensation of wage costs, up to 90% of wage bill, conditional on at least 20% fall in turnover in the months March to May compared to 2019 (EUR 10bn). 
- In addiiton to the wage scheme, the government will cover 30% of pension contributions and the employers' premium (1 April)
- Self-employed and independent entrepreneurs without staff can receive social assistance benefits for theree months through an accelerated procedure, where asset tests and tests on partner income are skipped. Loans at reduced interest rates are available under the same scheme (EUR 2bn).
- Self-employed (without employees) hit by the corona crisis are to receive monthly 1,050 'euros over the next three months, irrespective of their savings. For married couples or couples with children a maximum of EUR 1,500 a month will apply (30 March). 
- Ban on forced home sales of homeowners who are temporarily unable to meet their mortgage obligations due to the corona crisis until 1 July (7 April). Automatic extension of temporary rent contracts until the end of June (8 April)
- Childcare Compensation for parents who continue paying for childcare during the lockdown (EUR 175 Million) (8 April)

'@
$ws.Range("I38").Value = $text_I38

$text_J38 = @'
- Tax collection of VAT, income, corporate, wage and turnover taxes will be halted immediately upon application to defer tax payments from companies. The estimated size of this support is between EUR 45-60 bn (5.8%-8.4%) of GDP in 2020. Companies get more time to provide documentation. The interest rate on tax debts is reduced from 4% to almost 0%. On April 6, the tax holiday was extended to excise duties, environmental tax, and gambling and insurance tax (at least until 19 June) 
- A guarantee facility for SME loans (BMKB) is already in place and working – EUR 665 million of an available budget of EUR 765 million, 0.08% of GDP. The Ministry of Economic Affairs and Climate Policy has expanded guarantees to a higher guarantee percentage, up to 90% of the credit risk of banks for SMEs that want to take out a loan of EUR 1.5 mln maximum but that do not have the required collateral. On 9 April, budget for the MKB scheme was raised to EUR 1.5 billion and the premium for SME loans has been reduced from 3.9 to 2%.
- The Guarantee for Entrepreneurial Financing (GO-regeling) for businesses that experience problems in obtaining bank loans and bank guarantees has been broadened. The guarantee ceiling increases from 400 million to 1.5 billion euros. Both SMEs and large companies can obtain a 50% guarantee on bank loans and bank guarantees (minimum 1.5 million - maximum 50 million euros per company). The maximum per company is temporarily increased to 150 million euros. On 9 April, the budget for the Go-regeling scheme was raised to EUR 10 billion. 
-  Firms that needed to close due to the coronavurus can receive a one-off lump sum allowance of EUR 4000 (TOGS scheme) (30 March). In a first instance, this included firms in the hospitality sector, hairdressers and the travel industry. On 31 March, the emergency desk was extended to retailers. On 9 April, it was further extended to include tattoo shops, taxi companies, dentists, physiotherapists, and suppliers of the hospitality industry or events.
- The government provides EUR 12 billion guarantees for the credit insurance market to ease lending to small firms (9 April) 
- Separate credit option for start-ups and scale-ups, who can apply as Regional Development Companies (EUR 100 million) (9 April)
- The banking sector has responded by giving an aoutomatic 6-month payment holiday (interest and amortisation) for all business loans of less than 2.5 million euros. On 30 March, the bank ABN Amro raised the ceiling to EUR 50 million to extend the loan repayment pause to bigger enterprises. 
Interest rates are also lowered for a publicly backed microcredit scheme (Qredits).
- A temporary guarantee for working capital in agricultural and horticultural companies will be provided under the Guarantee for SME Agricultural Credits (BL) scheme.
- The government increased export credit insurance for firms through the credit insurer Atradius. The share of the working capital that companies need for their export production was raised from 80 to 95 percent (30 March)
- Competition regulations and driving and rest times regulations are eased temporarily to secure supply of foodstuffs, energy pharmaceutical products, as well as transport of waste and garbage.

'@
$ws.Range("J38").Value = $text_J38

$text_C49 = @'
• Mandated nation-wide quarantine for at least 15 days since March 15th, which has been further extended until April 11th and then to April 26th. 
• Suspension of all non-essential work from March 30 to April 9, inclusive.
'@
$ws.Range("C49").Value = $text_C49

